# racunice uz klizuci datum
# Append a new data row (row 10) to Sheet1 with a new station record:
#   ams=00000E98, opstina=Ruma, mesto=Irig-Kudos, pocetak=22.02.2011, usev=Jabuka

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "00000E98" looks numeric (scientific-notation-like) to Excel's auto-detection,
# so force it to stay text by entering it with a leading apostrophe, same as the
# other "ams" id codes in column A (which use the quote-prefixed text style).
$ws.Range("A10").Value = "'00000E98"
$ws.Range("B10").Value = "Ruma"
$ws.Range("C10").Value = "Irig-Kudos"
$ws.Range("D10").Value = "22.02.2011"
$ws.Range("F10").Value = "Jabuka"

[void]$ws.Range("E10").Select()
